$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Range("H11").Value = 135.6
$ws.Range("I11").Value = 135.6
$ws.Range("K11").Value = 135.6
$ws.Range("M11").Value = 4.400000000000006
# Row 18
$ws.Range("H18").Value = 1749.5
$ws.Range("I18").Value = 1999
$ws.Range("K18").Value = 1999
$ws.Range("M18").Value = -1715
# Row 80
$ws.Range("H80").Value = 702.5454999999999
$ws.Range("I80").Value = 828.1429000000001
$ws.Range("J80").Value = 482.75
$ws.Range("K80").Value = 2484.4287
$ws.Range("L80").Value = 1448.25
$ws.Range("M80").Value = -1486.4287
$ws.Range("N80").Value = -3444.25
# Row 83
$ws.Range("H83").Value = 702.5454999999999
$ws.Range("I83").Value = 828.1429000000001
$ws.Range("J83").Value = 482.75
$ws.Range("K83").Value = 7453.2861
$ws.Range("L83").Value = 4344.75
$ws.Range("M83").Value = -2461.2861
$ws.Range("N83").Value = -14328.75
# Row 98
$ws.Range("H98").Value = 1944.2142
$ws.Range("I98").Value = 1710.3
$ws.Range("K98").Value = 1710.3
$ws.Range("M98").Value = -212.3
# Row 121
$ws.Range("H121").Value = 2395.077
$ws.Range("J121").Value = 2395.077
$ws.Range("L121").Value = 7185.231000000001
$ws.Range("N121").Value = -10679.231
# Row 122
$ws.Range("H122").Value = 1944.2142
$ws.Range("I122").Value = 1710.3
$ws.Range("K122").Value = 5130.9
$ws.Range("M122").Value = -2680.9
# Row 131
$ws.Range("H131").Value = 1918.4117
$ws.Range("I131").Value = 1538.6875
$ws.Range("J131").Value = 7994
$ws.Range("K131").Value = 4616.0625
$ws.Range("L131").Value = 23982
$ws.Range("M131").Value = 423.9375
$ws.Range("N131").Value = -34062
# Row 132
$ws.Range("H132").Value = 8719.023999999999
$ws.Range("I132").Value = 6409.59
$ws.Range("K132").Value = 19228.77
$ws.Range("M132").Value = -16698.77
# Row 137
$ws.Range("H137").Value = 10755.489
$ws.Range("I137").Value = 3319.4644
$ws.Range("J137").Value = 23003.059
$ws.Range("K137").Value = 9958.393199999999
$ws.Range("L137").Value = 69009.177
$ws.Range("M137").Value = -7408.393199999999
$ws.Range("N137").Value = -74109.177

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 12278.758
$ws.Range("I32").Value = 4477.551
$ws.Range("K32").Value = 4477.551
$ws.Range("M32").Value = -4190.551
# Row 74
$ws.Range("H74").Value = 14880.4375
$ws.Range("I74").Value = 2178.4285
$ws.Range("K74").Value = 2178.4285
$ws.Range("M74").Value = -1304.4285
# Row 77
$ws.Range("H77").Value = 14880.4375
$ws.Range("I77").Value = 2178.4285
$ws.Range("K77").Value = 10892.1425
$ws.Range("M77").Value = -6524.1425
# Row 138
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 4212.6
$ws.Range("I86").Value = 3823.625
$ws.Range("J86").Value = 4657.143
$ws.Range("K86").Value = 3823.625
$ws.Range("L86").Value = 4657.143
$ws.Range("M86").Value = -2700.625
$ws.Range("N86").Value = -6903.143
# Row 89
$ws.Range("H89").Value = 4212.6
$ws.Range("I89").Value = 3823.625
$ws.Range("J89").Value = 4657.143
$ws.Range("K89").Value = 19118.125
$ws.Range("L89").Value = 23285.715
$ws.Range("M89").Value = -13502.125
$ws.Range("N89").Value = -34517.715
# Row 105
$ws.Range("H105").Value = 2080
$ws.Range("I105").Value = 1854.1177
$ws.Range("K105").Value = 1854.1177
$ws.Range("M105").Value = -107.1177
# Row 134
$ws.Range("H134").Value = 8960.093000000001
$ws.Range("I134").Value = 3977.5151
$ws.Range("K134").Value = 11932.5453
$ws.Range("M134").Value = -9397.5453

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 21453.432
$ws.Range("I31").Value = 12601.8
$ws.Range("J31").Value = 24731.814
$ws.Range("K31").Value = 12601.8
$ws.Range("L31").Value = 24731.814
$ws.Range("M31").Value = -12306.8
$ws.Range("N31").Value = -25321.814
# Row 34
$ws.Range("H34").Value = 21453.432
$ws.Range("I34").Value = 12601.8
$ws.Range("J34").Value = 24731.814
$ws.Range("K34").Value = 12601.8
$ws.Range("L34").Value = 24731.814
$ws.Range("M34").Value = -12399.8
$ws.Range("N34").Value = -25135.814
# Row 41
$ws.Range("H41").Value = 17084.8
$ws.Range("I41").Value = 12529.5
$ws.Range("J41").Value = 20121.666
$ws.Range("K41").Value = 12529.5
$ws.Range("L41").Value = 20121.666
$ws.Range("M41").Value = -12101.5
$ws.Range("N41").Value = -20977.666
# Row 50
$ws.Range("H50").Value = 41333
$ws.Range("J50").Value = 41333
$ws.Range("L50").Value = 41333
$ws.Range("N50").Value = -42583
# Row 58
$ws.Range("H58").Value = 13512.522
$ws.Range("I58").Value = 5620.0625
$ws.Range("K58").Value = 5620.0625
$ws.Range("M58").Value = -5417.0625
# Row 59
$ws.Range("H59").Value = 44999
$ws.Range("J59").Value = 44999
$ws.Range("L59").Value = 44999
$ws.Range("N59").Value = -47289
# Row 60
$ws.Range("H60").Value = 28398.6
$ws.Range("J60").Value = 33750
$ws.Range("L60").Value = 33750
$ws.Range("N60").Value = -34772
# Row 62
$ws.Range("H62").Value = 4874.5
$ws.Range("I62").Value = 4049.5
$ws.Range("J62").Value = 7349.5
$ws.Range("K62").Value = 4049.5
$ws.Range("L62").Value = 7349.5
$ws.Range("M62").Value = -3425.5
$ws.Range("N62").Value = -8597.5
# Row 65
$ws.Range("H65").Value = 4874.5
$ws.Range("I65").Value = 4049.5
$ws.Range("J65").Value = 7349.5
$ws.Range("K65").Value = 20247.5
$ws.Range("L65").Value = 36747.5
$ws.Range("M65").Value = -17127.5
$ws.Range("N65").Value = -42987.5
# Row 86
$ws.Range("H86").Value = 3705.5
$ws.Range("I86").Value = 2570.524
$ws.Range("J86").Value = 6353.778
$ws.Range("K86").Value = 2570.524
$ws.Range("L86").Value = 6353.778
$ws.Range("M86").Value = -1447.524
$ws.Range("N86").Value = -8599.778
# Row 89
$ws.Range("H89").Value = 3705.5
$ws.Range("I89").Value = 2570.524
$ws.Range("J89").Value = 6353.778
$ws.Range("K89").Value = 12852.62
$ws.Range("L89").Value = 31768.89
$ws.Range("M89").Value = -7236.619999999999
$ws.Range("N89").Value = -43000.89
# Row 94
$ws.Range("H94").Value = 1799.8
$ws.Range("I94").Value = 2002.6666
$ws.Range("J94").Value = 1495.5
$ws.Range("K94").Value = 2002.6666
$ws.Range("L94").Value = 1495.5
$ws.Range("M94").Value = -1551.6666
$ws.Range("N94").Value = -2397.5
# Row 132
$ws.Range("H132").Value = 6803.625
$ws.Range("I132").Value = 2606.4
$ws.Range("J132").Value = 10507.059
$ws.Range("K132").Value = 7819.200000000001
$ws.Range("L132").Value = 31521.177
$ws.Range("M132").Value = -5289.200000000001
$ws.Range("N132").Value = -36581.177
# Row 134
$ws.Range("H134").Value = 20412348
$ws.Range("I134").Value = 1069.25
$ws.Range("K134").Value = 3207.75
$ws.Range("M134").Value = -672.75
# Row 136
$ws.Range("H136").Value = 13512.522
$ws.Range("I136").Value = 5620.0625
$ws.Range("K136").Value = 16860.1875
$ws.Range("M136").Value = -14310.1875
# Row 141
$ws.Range("H141").Value = 139450.78
$ws.Range("J141").Value = 152988.75
$ws.Range("L141").Value = 152988.75
$ws.Range("N141").Value = -163348.75

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 63
$ws.Range("H63").Value = 10786.3
$ws.Range("I63").Value = 8644
$ws.Range("J63").Value = 13999.75
$ws.Range("K63").Value = 25932
$ws.Range("L63").Value = 41999.25
$ws.Range("M63").Value = -25183
$ws.Range("N63").Value = -43497.25
# Row 66
$ws.Range("H66").Value = 10786.3
$ws.Range("I66").Value = 8644
$ws.Range("J66").Value = 13999.75
$ws.Range("K66").Value = 77796
$ws.Range("L66").Value = 125997.75
$ws.Range("M66").Value = -74052
$ws.Range("N66").Value = -133485.75
# Row 138
$ws.Range("H138").Value = 3529.9722
$ws.Range("I138").Value = 931.9
$ws.Range("J138").Value = 4529.231
$ws.Range("K138").Value = 2795.7
$ws.Range("L138").Value = 13587.693
$ws.Range("M138").Value = 2344.3
$ws.Range("N138").Value = -23867.693

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 36
$ws.Range("H36").Value = 19486.875
$ws.Range("J36").Value = 20127.857
$ws.Range("L36").Value = 20127.857
$ws.Range("N36").Value = -21097.857
# Row 63
$ws.Range("H63").Value = 40000
$ws.Range("J63").Value = 30000
$ws.Range("L63").Value = 30000
$ws.Range("N63").Value = -31372
# Row 66
$ws.Range("H66").Value = 40000
$ws.Range("J66").Value = 30000
$ws.Range("L66").Value = 90000
$ws.Range("N66").Value = -96864

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 100
$ws.Range("H100").Value = 14874.75
$ws.Range("I100").Value = 3166.6667
$ws.Range("J100").Value = 49999
$ws.Range("K100").Value = 3166.6667
$ws.Range("L100").Value = 49999
$ws.Range("M100").Value = -2625.6667
$ws.Range("N100").Value = -51081
